$wb = $excel.ActiveWorkbook

# F-column ("想去人数") updates that apply identically to both the
# "展览" sheet and the aggregated "全部类型" sheet.
$updates = @{
    "F2"  = 1032
    "F3"  = 307
    "F4"  = 1421
    "F5"  = 8561
    "F6"  = 68
    "F9"  = 257
    "F11" = 3454
    "F13" = 347
    "F14" = 68
    "F15" = 1010
    "F17" = 1099
    "F18" = 300
    "F19" = 174
    "F20" = 2138
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
